$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$ws.Range("H3").Value = 286
$ws.Range("I3").Value = "'04-Nov-2025"

$ws.Range("H4").Value = 377
$ws.Range("I4").Value = "'04-Nov-2025"
